$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Notified Production (MW)" values for rows 2-97 (column B)
$bValues = @(915,920,917,918,973,980,982,982,1004,1003,1003,1004,1106,1108,1111,1114,1138,1130,1138,1131,1190,1187,1185,1185,1179,1159,1154,1153,1039,1030,1022,1012,855,852,849,846,933,934,935,935,922,922,923,923,907,906,904,903,909,907,906,904,831,829,827,825,722,720,717,714,559,555,552,549,379,377,375,374,312,313,314,315,283,283,284,284,267,267,267,267,257,257,258,258,227,227,227,228,199,199,199,199,0,0,0,0)

for ($i = 0; $i -lt $bValues.Length; $i++) {
    $row = $i + 2
    # Shift the timestamp in column A forward by exactly one day
    $oldDate = $ws.Cells.Item($row, 1).Value2
    $ws.Cells.Item($row, 1).Value = $oldDate + 1
    # Write the updated production value in column B
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}

Write-Output "Updated $($bValues.Length) rows"
